$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing Acelga price-history rows (182-301): a new weekly
# observation is inserted at the top and the series shifts down,
# appending the two oldest rows (300-301) at the bottom.

$dArr = New-Object 'object[,]' 120,1
$dArr[0,0] = 44452
$dArr[1,0] = 44452
$dArr[2,0] = 44358
$dArr[3,0] = 44358
$dArr[4,0] = 44358
$dArr[5,0] = 44218
$dArr[6,0] = 44218
$dArr[7,0] = 44433
$dArr[8,0] = 44433
$dArr[9,0] = 44433
$dArr[10,0] = 44397
$dArr[11,0] = 44397
$dArr[12,0] = 44397
$dArr[13,0] = 44274
$dArr[14,0] = 44274
$dArr[15,0] = 44426
$dArr[16,0] = 44426
$dArr[17,0] = 44421
$dArr[18,0] = 44421
$dArr[19,0] = 44421
$dArr[20,0] = 44434
$dArr[21,0] = 44434
$dArr[22,0] = 44434
$dArr[23,0] = 44215
$dArr[24,0] = 44194
$dArr[25,0] = 44194
$dArr[26,0] = 44420
$dArr[27,0] = 44420
$dArr[28,0] = 44420
$dArr[29,0] = 44264
$dArr[30,0] = 44264
$dArr[31,0] = 44348
$dArr[32,0] = 44348
$dArr[33,0] = 44348
$dArr[34,0] = 44316
$dArr[35,0] = 44316
$dArr[36,0] = 44253
$dArr[37,0] = 44253
$dArr[38,0] = 44427
$dArr[39,0] = 44427
$dArr[40,0] = 44427
$dArr[41,0] = 44341
$dArr[42,0] = 44341
$dArr[43,0] = 44341
$dArr[44,0] = 44414
$dArr[45,0] = 44414
$dArr[46,0] = 44414
$dArr[47,0] = 44301
$dArr[48,0] = 44301
$dArr[49,0] = 44176
$dArr[50,0] = 44441
$dArr[51,0] = 44441
$dArr[52,0] = 44441
$dArr[53,0] = 44432
$dArr[54,0] = 44432
$dArr[55,0] = 44432
$dArr[56,0] = 44351
$dArr[57,0] = 44351
$dArr[58,0] = 44351
$dArr[59,0] = 44369
$dArr[60,0] = 44369
$dArr[61,0] = 44369
$dArr[62,0] = 44257
$dArr[63,0] = 44257
$dArr[64,0] = 44273
$dArr[65,0] = 44273
$dArr[66,0] = 44294
$dArr[67,0] = 44294
$dArr[68,0] = 44379
$dArr[69,0] = 44379
$dArr[70,0] = 44379
$dArr[71,0] = 44302
$dArr[72,0] = 44302
$dArr[73,0] = 44315
$dArr[74,0] = 44315
$dArr[75,0] = 44391
$dArr[76,0] = 44391
$dArr[77,0] = 44446
$dArr[78,0] = 44446
$dArr[79,0] = 44446
$dArr[80,0] = 44411
$dArr[81,0] = 44411
$dArr[82,0] = 44411
$dArr[83,0] = 44313
$dArr[84,0] = 44313
$dArr[85,0] = 44329
$dArr[86,0] = 44329
$dArr[87,0] = 44329
$dArr[88,0] = 44449
$dArr[89,0] = 44449
$dArr[90,0] = 44449
$dArr[91,0] = 44161
$dArr[92,0] = 44161
$dArr[93,0] = 44428
$dArr[94,0] = 44428
$dArr[95,0] = 44428
$dArr[96,0] = 44442
$dArr[97,0] = 44442
$dArr[98,0] = 44442
$dArr[99,0] = 44435
$dArr[100,0] = 44435
$dArr[101,0] = 44435
$dArr[102,0] = 44319
$dArr[103,0] = 44175
$dArr[104,0] = 44376
$dArr[105,0] = 44376
$dArr[106,0] = 44376
$dArr[107,0] = 44223
$dArr[108,0] = 44223
$dArr[109,0] = 44448
$dArr[110,0] = 44448
$dArr[111,0] = 44448
$dArr[112,0] = 44238
$dArr[113,0] = 44238
$dArr[114,0] = 44399
$dArr[115,0] = 44399
$dArr[116,0] = 44399
$dArr[117,0] = 44400
$dArr[118,0] = 44400
$dArr[119,0] = 44400
$ws.Range("D182:D301").Value = $dArr

$imArr = New-Object 'object[,]' 120,5
$imArr[0,0] = "Primera"
$imArr[0,1] = 34
$imArr[0,2] = 11000
$imArr[0,3] = 12000
$imArr[0,4] = 11500
$imArr[1,0] = "Segunda"
$imArr[1,1] = 16
$imArr[1,2] = 9000
$imArr[1,3] = 10000
$imArr[1,4] = 9500
$imArr[2,0] = "Extra"
$imArr[2,1] = 16
$imArr[2,2] = 13000
$imArr[2,3] = 13000
$imArr[2,4] = 13000
$imArr[3,0] = "Primera"
$imArr[3,1] = 52
$imArr[3,2] = 10000
$imArr[3,3] = 11000
$imArr[3,4] = 10500
$imArr[4,0] = "Segunda"
$imArr[4,1] = 34
$imArr[4,2] = 9000
$imArr[4,3] = 9000
$imArr[4,4] = 9000
$imArr[5,0] = "Primera"
$imArr[5,1] = 80
$imArr[5,2] = 10000
$imArr[5,3] = 12000
$imArr[5,4] = 11000
$imArr[6,0] = "Segunda"
$imArr[6,1] = 45
$imArr[6,2] = 8000
$imArr[6,3] = 8000
$imArr[6,4] = 8000
$imArr[7,0] = "Extra"
$imArr[7,1] = 16
$imArr[7,2] = 13000
$imArr[7,3] = 13000
$imArr[7,4] = 13000
$imArr[8,0] = "Primera"
$imArr[8,1] = 43
$imArr[8,2] = 11000
$imArr[8,3] = 12000
$imArr[8,4] = 11512
$imArr[9,0] = "Segunda"
$imArr[9,1] = 25
$imArr[9,2] = 9000
$imArr[9,3] = 10000
$imArr[9,4] = 9480
$imArr[10,0] = "Extra"
$imArr[10,1] = 16
$imArr[10,2] = 13000
$imArr[10,3] = 13000
$imArr[10,4] = 13000
$imArr[11,0] = "Primera"
$imArr[11,1] = 52
$imArr[11,2] = 10000
$imArr[11,3] = 11000
$imArr[11,4] = 10500
$imArr[12,0] = "Segunda"
$imArr[12,1] = 25
$imArr[12,2] = 9000
$imArr[12,3] = 9000
$imArr[12,4] = 9000
$imArr[13,0] = "Primera"
$imArr[13,1] = 70
$imArr[13,2] = 15000
$imArr[13,3] = 15000
$imArr[13,4] = 15000
$imArr[14,0] = "Segunda"
$imArr[14,1] = 50
$imArr[14,2] = 13000
$imArr[14,3] = 13000
$imArr[14,4] = 13000
$imArr[15,0] = "Primera"
$imArr[15,1] = 52
$imArr[15,2] = 11000
$imArr[15,3] = 12000
$imArr[15,4] = 11500
$imArr[16,0] = "Segunda"
$imArr[16,1] = 25
$imArr[16,2] = 10000
$imArr[16,3] = 10000
$imArr[16,4] = 10000
$imArr[17,0] = "Extra"
$imArr[17,1] = 25
$imArr[17,2] = 13000
$imArr[17,3] = 13000
$imArr[17,4] = 13000
$imArr[18,0] = "Primera"
$imArr[18,1] = 61
$imArr[18,2] = 11000
$imArr[18,3] = 12000
$imArr[18,4] = 11492
$imArr[19,0] = "Segunda"
$imArr[19,1] = 43
$imArr[19,2] = 9000
$imArr[19,3] = 10000
$imArr[19,4] = 9488
$imArr[20,0] = "Extra"
$imArr[20,1] = 25
$imArr[20,2] = 14000
$imArr[20,3] = 14000
$imArr[20,4] = 14000
$imArr[21,0] = "Primera"
$imArr[21,1] = 43
$imArr[21,2] = 12000
$imArr[21,3] = 13000
$imArr[21,4] = 12512
$imArr[22,0] = "Segunda"
$imArr[22,1] = 34
$imArr[22,2] = 10000
$imArr[22,3] = 11000
$imArr[22,4] = 10500
$imArr[23,0] = "Primera"
$imArr[23,1] = 65
$imArr[23,2] = 10000
$imArr[23,3] = 12000
$imArr[23,4] = 10769
$imArr[24,0] = "Primera"
$imArr[24,1] = 70
$imArr[24,2] = 10000
$imArr[24,3] = 10000
$imArr[24,4] = 10000
$imArr[25,0] = "Segunda"
$imArr[25,1] = 30
$imArr[25,2] = 8000
$imArr[25,3] = 8000
$imArr[25,4] = 8000
$imArr[26,0] = "Extra"
$imArr[26,1] = 25
$imArr[26,2] = 13000
$imArr[26,3] = 13000
$imArr[26,4] = 13000
$imArr[27,0] = "Primera"
$imArr[27,1] = 61
$imArr[27,2] = 11000
$imArr[27,3] = 12000
$imArr[27,4] = 11492
$imArr[28,0] = "Segunda"
$imArr[28,1] = 34
$imArr[28,2] = 10000
$imArr[28,3] = 10000
$imArr[28,4] = 10000
$imArr[29,0] = "Primera"
$imArr[29,1] = 70
$imArr[29,2] = 14000
$imArr[29,3] = 14000
$imArr[29,4] = 14000
$imArr[30,0] = "Segunda"
$imArr[30,1] = 50
$imArr[30,2] = 12000
$imArr[30,3] = 12000
$imArr[30,4] = 12000
$imArr[31,0] = "Extra"
$imArr[31,1] = 20
$imArr[31,2] = 14000
$imArr[31,3] = 14000
$imArr[31,4] = 14000
$imArr[32,0] = "Primera"
$imArr[32,1] = 50
$imArr[32,2] = 11000
$imArr[32,3] = 12000
$imArr[32,4] = 11500
$imArr[33,0] = "Segunda"
$imArr[33,1] = 30
$imArr[33,2] = 9000
$imArr[33,3] = 10000
$imArr[33,4] = 9500
$imArr[34,0] = "Primera"
$imArr[34,1] = 30
$imArr[34,2] = 13000
$imArr[34,3] = 13000
$imArr[34,4] = 13000
$imArr[35,0] = "Segunda"
$imArr[35,1] = 16
$imArr[35,2] = 10000
$imArr[35,3] = 10000
$imArr[35,4] = 10000
$imArr[36,0] = "Primera"
$imArr[36,1] = 70
$imArr[36,2] = 15000
$imArr[36,3] = 15000
$imArr[36,4] = 15000
$imArr[37,0] = "Segunda"
$imArr[37,1] = 50
$imArr[37,2] = 12000
$imArr[37,3] = 12000
$imArr[37,4] = 12000
$imArr[38,0] = "Extra"
$imArr[38,1] = 16
$imArr[38,2] = 14000
$imArr[38,3] = 14000
$imArr[38,4] = 14000
$imArr[39,0] = "Primera"
$imArr[39,1] = 34
$imArr[39,2] = 12000
$imArr[39,3] = 13000
$imArr[39,4] = 12500
$imArr[40,0] = "Segunda"
$imArr[40,1] = 25
$imArr[40,2] = 10000
$imArr[40,3] = 11000
$imArr[40,4] = 10520
$imArr[41,0] = "Extra"
$imArr[41,1] = 16
$imArr[41,2] = 14000
$imArr[41,3] = 14000
$imArr[41,4] = 14000
$imArr[42,0] = "Primera"
$imArr[42,1] = 50
$imArr[42,2] = 12000
$imArr[42,3] = 12000
$imArr[42,4] = 12000
$imArr[43,0] = "Segunda"
$imArr[43,1] = 13
$imArr[43,2] = 10000
$imArr[43,3] = 10000
$imArr[43,4] = 10000
$imArr[44,0] = "Extra"
$imArr[44,1] = 25
$imArr[44,2] = 13000
$imArr[44,3] = 14000
$imArr[44,4] = 13480
$imArr[45,0] = "Primera"
$imArr[45,1] = 70
$imArr[45,2] = 11000
$imArr[45,3] = 12000
$imArr[45,4] = 11500
$imArr[46,0] = "Segunda"
$imArr[46,1] = 34
$imArr[46,2] = 9000
$imArr[46,3] = 10000
$imArr[46,4] = 9500
$imArr[47,0] = "Primera"
$imArr[47,1] = 50
$imArr[47,2] = 15000
$imArr[47,3] = 15000
$imArr[47,4] = 15000
$imArr[48,0] = "Segunda"
$imArr[48,1] = 30
$imArr[48,2] = 13000
$imArr[48,3] = 13000
$imArr[48,4] = 13000
$imArr[49,0] = "Primera"
$imArr[49,1] = 250
$imArr[49,2] = 10000
$imArr[49,3] = 11000
$imArr[49,4] = 10500
$imArr[50,0] = "Extra"
$imArr[50,1] = 16
$imArr[50,2] = 14000
$imArr[50,3] = 14000
$imArr[50,4] = 14000
$imArr[51,0] = "Primera"
$imArr[51,1] = 34
$imArr[51,2] = 12000
$imArr[51,3] = 13000
$imArr[51,4] = 12500
$imArr[52,0] = "Segunda"
$imArr[52,1] = 25
$imArr[52,2] = 10000
$imArr[52,3] = 11000
$imArr[52,4] = 10480
$imArr[53,0] = "Extra"
$imArr[53,1] = 16
$imArr[53,2] = 14000
$imArr[53,3] = 14000
$imArr[53,4] = 14000
$imArr[54,0] = "Primera"
$imArr[54,1] = 52
$imArr[54,2] = 12000
$imArr[54,3] = 13000
$imArr[54,4] = 12500
$imArr[55,0] = "Segunda"
$imArr[55,1] = 34
$imArr[55,2] = 10000
$imArr[55,3] = 11000
$imArr[55,4] = 10500
$imArr[56,0] = "Extra"
$imArr[56,1] = 25
$imArr[56,2] = 14000
$imArr[56,3] = 14000
$imArr[56,4] = 14000
$imArr[57,0] = "Primera"
$imArr[57,1] = 52
$imArr[57,2] = 11000
$imArr[57,3] = 12000
$imArr[57,4] = 11500
$imArr[58,0] = "Segunda"
$imArr[58,1] = 30
$imArr[58,2] = 9000
$imArr[58,3] = 9000
$imArr[58,4] = 9000
$imArr[59,0] = "Extra"
$imArr[59,1] = 25
$imArr[59,2] = 14000
$imArr[59,3] = 14000
$imArr[59,4] = 14000
$imArr[60,0] = "Primera"
$imArr[60,1] = 52
$imArr[60,2] = 12000
$imArr[60,3] = 13000
$imArr[60,4] = 12500
$imArr[61,0] = "Segunda"
$imArr[61,1] = 34
$imArr[61,2] = 10000
$imArr[61,3] = 10000
$imArr[61,4] = 10000
$imArr[62,0] = "Primera"
$imArr[62,1] = 70
$imArr[62,2] = 15000
$imArr[62,3] = 15000
$imArr[62,4] = 15000
$imArr[63,0] = "Segunda"
$imArr[63,1] = 50
$imArr[63,2] = 12000
$imArr[63,3] = 12000
$imArr[63,4] = 12000
$imArr[64,0] = "Primera"
$imArr[64,1] = 70
$imArr[64,2] = 14000
$imArr[64,3] = 14000
$imArr[64,4] = 14000
$imArr[65,0] = "Segunda"
$imArr[65,1] = 50
$imArr[65,2] = 12000
$imArr[65,3] = 12000
$imArr[65,4] = 12000
$imArr[66,0] = "Primera"
$imArr[66,1] = 70
$imArr[66,2] = 15000
$imArr[66,3] = 15000
$imArr[66,4] = 15000
$imArr[67,0] = "Segunda"
$imArr[67,1] = 50
$imArr[67,2] = 12000
$imArr[67,3] = 12000
$imArr[67,4] = 12000
$imArr[68,0] = "Extra"
$imArr[68,1] = 16
$imArr[68,2] = 13000
$imArr[68,3] = 13000
$imArr[68,4] = 13000
$imArr[69,0] = "Primera"
$imArr[69,1] = 52
$imArr[69,2] = 11000
$imArr[69,3] = 12000
$imArr[69,4] = 11500
$imArr[70,0] = "Segunda"
$imArr[70,1] = 34
$imArr[70,2] = 9000
$imArr[70,3] = 9000
$imArr[70,4] = 9000
$imArr[71,0] = "Primera"
$imArr[71,1] = 50
$imArr[71,2] = 15000
$imArr[71,3] = 15000
$imArr[71,4] = 15000
$imArr[72,0] = "Segunda"
$imArr[72,1] = 30
$imArr[72,2] = 13000
$imArr[72,3] = 13000
$imArr[72,4] = 13000
$imArr[73,0] = "Primera"
$imArr[73,1] = 40
$imArr[73,2] = 12000
$imArr[73,3] = 12000
$imArr[73,4] = 12000
$imArr[74,0] = "Segunda"
$imArr[74,1] = 20
$imArr[74,2] = 10000
$imArr[74,3] = 10000
$imArr[74,4] = 10000
$imArr[75,0] = "Primera"
$imArr[75,1] = 70
$imArr[75,2] = 11000
$imArr[75,3] = 12000
$imArr[75,4] = 11500
$imArr[76,0] = "Segunda"
$imArr[76,1] = 43
$imArr[76,2] = 9000
$imArr[76,3] = 9000
$imArr[76,4] = 9000
$imArr[77,0] = "Extra"
$imArr[77,1] = 16
$imArr[77,2] = 13000
$imArr[77,3] = 13000
$imArr[77,4] = 13000
$imArr[78,0] = "Primera"
$imArr[78,1] = 52
$imArr[78,2] = 11000
$imArr[78,3] = 12000
$imArr[78,4] = 11500
$imArr[79,0] = "Segunda"
$imArr[79,1] = 34
$imArr[79,2] = 9000
$imArr[79,3] = 10000
$imArr[79,4] = 9500
$imArr[80,0] = "Extra"
$imArr[80,1] = 16
$imArr[80,2] = 15000
$imArr[80,3] = 15000
$imArr[80,4] = 15000
$imArr[81,0] = "Primera"
$imArr[81,1] = 52
$imArr[81,2] = 12000
$imArr[81,3] = 13000
$imArr[81,4] = 12500
$imArr[82,0] = "Segunda"
$imArr[82,1] = 34
$imArr[82,2] = 11000
$imArr[82,3] = 11000
$imArr[82,4] = 11000
$imArr[83,0] = "Primera"
$imArr[83,1] = 70
$imArr[83,2] = 15000
$imArr[83,3] = 15000
$imArr[83,4] = 15000
$imArr[84,0] = "Segunda"
$imArr[84,1] = 20
$imArr[84,2] = 13000
$imArr[84,3] = 13000
$imArr[84,4] = 13000
$imArr[85,0] = "Extra"
$imArr[85,1] = 25
$imArr[85,2] = 15000
$imArr[85,3] = 15000
$imArr[85,4] = 15000
$imArr[86,0] = "Primera"
$imArr[86,1] = 50
$imArr[86,2] = 13000
$imArr[86,3] = 13000
$imArr[86,4] = 13000
$imArr[87,0] = "Segunda"
$imArr[87,1] = 20
$imArr[87,2] = 11000
$imArr[87,3] = 11000
$imArr[87,4] = 11000
$imArr[88,0] = "Extra"
$imArr[88,1] = 25
$imArr[88,2] = 12000
$imArr[88,3] = 13000
$imArr[88,4] = 12520
$imArr[89,0] = "Primera"
$imArr[89,1] = 52
$imArr[89,2] = 10000
$imArr[89,3] = 11000
$imArr[89,4] = 10500
$imArr[90,0] = "Segunda"
$imArr[90,1] = 34
$imArr[90,2] = 8000
$imArr[90,3] = 9000
$imArr[90,4] = 8500
$imArr[91,0] = "Primera"
$imArr[91,1] = 75
$imArr[91,2] = 10000
$imArr[91,3] = 12000
$imArr[91,4] = 11067
$imArr[92,0] = "Segunda"
$imArr[92,1] = 50
$imArr[92,2] = 8000
$imArr[92,3] = 8000
$imArr[92,4] = 8000
$imArr[93,0] = "Extra"
$imArr[93,1] = 16
$imArr[93,2] = 14000
$imArr[93,3] = 14000
$imArr[93,4] = 14000
$imArr[94,0] = "Primera"
$imArr[94,1] = 34
$imArr[94,2] = 12000
$imArr[94,3] = 13000
$imArr[94,4] = 12500
$imArr[95,0] = "Segunda"
$imArr[95,1] = 25
$imArr[95,2] = 10000
$imArr[95,3] = 11000
$imArr[95,4] = 10520
$imArr[96,0] = "Extra"
$imArr[96,1] = 18
$imArr[96,2] = 13000
$imArr[96,3] = 13000
$imArr[96,4] = 13000
$imArr[97,0] = "Primera"
$imArr[97,1] = 38
$imArr[97,2] = 11000
$imArr[97,3] = 12000
$imArr[97,4] = 11500
$imArr[98,0] = "Segunda"
$imArr[98,1] = 28
$imArr[98,2] = 9000
$imArr[98,3] = 10000
$imArr[98,4] = 9500
$imArr[99,0] = "Extra"
$imArr[99,1] = 70
$imArr[99,2] = 13000
$imArr[99,3] = 14000
$imArr[99,4] = 13586
$imArr[100,0] = "Primera"
$imArr[100,1] = 170
$imArr[100,2] = 11000
$imArr[100,3] = 13000
$imArr[100,4] = 12065
$imArr[101,0] = "Segunda"
$imArr[101,1] = 116
$imArr[101,2] = 9000
$imArr[101,3] = 11000
$imArr[101,4] = 10086
$imArr[102,0] = "Primera"
$imArr[102,1] = 34
$imArr[102,2] = 16000
$imArr[102,3] = 16000
$imArr[102,4] = 16000
$imArr[103,0] = "Primera"
$imArr[103,1] = 160
$imArr[103,2] = 10000
$imArr[103,3] = 11000
$imArr[103,4] = 10500
$imArr[104,0] = "Extra"
$imArr[104,1] = 16
$imArr[104,2] = 14000
$imArr[104,3] = 14000
$imArr[104,4] = 14000
$imArr[105,0] = "Primera"
$imArr[105,1] = 34
$imArr[105,2] = 11000
$imArr[105,3] = 12000
$imArr[105,4] = 11500
$imArr[106,0] = "Segunda"
$imArr[106,1] = 25
$imArr[106,2] = 9000
$imArr[106,3] = 9000
$imArr[106,4] = 9000
$imArr[107,0] = "Primera"
$imArr[107,1] = 60
$imArr[107,2] = 10000
$imArr[107,3] = 12000
$imArr[107,4] = 11333
$imArr[108,0] = "Segunda"
$imArr[108,1] = 40
$imArr[108,2] = 8000
$imArr[108,3] = 8000
$imArr[108,4] = 8000
$imArr[109,0] = "Extra"
$imArr[109,1] = 16
$imArr[109,2] = 12000
$imArr[109,3] = 13000
$imArr[109,4] = 12500
$imArr[110,0] = "Primera"
$imArr[110,1] = 43
$imArr[110,2] = 10000
$imArr[110,3] = 11000
$imArr[110,4] = 10488
$imArr[111,0] = "Segunda"
$imArr[111,1] = 34
$imArr[111,2] = 8000
$imArr[111,3] = 9000
$imArr[111,4] = 8500
$imArr[112,0] = "Primera"
$imArr[112,1] = 50
$imArr[112,2] = 14000
$imArr[112,3] = 14000
$imArr[112,4] = 14000
$imArr[113,0] = "Segunda"
$imArr[113,1] = 30
$imArr[113,2] = 10000
$imArr[113,3] = 10000
$imArr[113,4] = 10000
$imArr[114,0] = "Extra"
$imArr[114,1] = 18
$imArr[114,2] = 12000
$imArr[114,3] = 12000
$imArr[114,4] = 12000
$imArr[115,0] = "Primera"
$imArr[115,1] = 34
$imArr[115,2] = 10000
$imArr[115,3] = 11000
$imArr[115,4] = 10500
$imArr[116,0] = "Segunda"
$imArr[116,1] = 25
$imArr[116,2] = 9000
$imArr[116,3] = 9000
$imArr[116,4] = 9000
$imArr[117,0] = "Extra"
$imArr[117,1] = 16
$imArr[117,2] = 12000
$imArr[117,3] = 12000
$imArr[117,4] = 12000
$imArr[118,0] = "Primera"
$imArr[118,1] = 52
$imArr[118,2] = 10000
$imArr[118,3] = 11000
$imArr[118,4] = 10500
$imArr[119,0] = "Segunda"
$imArr[119,1] = 34
$imArr[119,2] = 8000
$imArr[119,3] = 9000
$imArr[119,4] = 8500
$ws.Range("I182:M301").Value = $imArr

$pArr = New-Object 'object[,]' 120,1
$pArr[0,0] = 3833
$pArr[1,0] = 3167
$pArr[2,0] = 4333
$pArr[3,0] = 3500
$pArr[4,0] = 3000
$pArr[5,0] = 3667
$pArr[6,0] = 2667
$pArr[7,0] = 4333
$pArr[8,0] = 3837
$pArr[9,0] = 3160
$pArr[10,0] = 4333
$pArr[11,0] = 3500
$pArr[12,0] = 3000
$pArr[13,0] = 5000
$pArr[14,0] = 4333
$pArr[15,0] = 3833
$pArr[16,0] = 3333
$pArr[17,0] = 4333
$pArr[18,0] = 3831
$pArr[19,0] = 3163
$pArr[20,0] = 4667
$pArr[21,0] = 4171
$pArr[22,0] = 3500
$pArr[23,0] = 3590
$pArr[24,0] = 3333
$pArr[25,0] = 2667
$pArr[26,0] = 4333
$pArr[27,0] = 3831
$pArr[28,0] = 3333
$pArr[29,0] = 4667
$pArr[30,0] = 4000
$pArr[31,0] = 4667
$pArr[32,0] = 3833
$pArr[33,0] = 3167
$pArr[34,0] = 4333
$pArr[35,0] = 3333
$pArr[36,0] = 5000
$pArr[37,0] = 4000
$pArr[38,0] = 4667
$pArr[39,0] = 4167
$pArr[40,0] = 3507
$pArr[41,0] = 4667
$pArr[42,0] = 4000
$pArr[43,0] = 3333
$pArr[44,0] = 4493
$pArr[45,0] = 3833
$pArr[46,0] = 3167
$pArr[47,0] = 5000
$pArr[48,0] = 4333
$pArr[49,0] = 3500
$pArr[50,0] = 4667
$pArr[51,0] = 4167
$pArr[52,0] = 3493
$pArr[53,0] = 4667
$pArr[54,0] = 4167
$pArr[55,0] = 3500
$pArr[56,0] = 4667
$pArr[57,0] = 3833
$pArr[58,0] = 3000
$pArr[59,0] = 4667
$pArr[60,0] = 4167
$pArr[61,0] = 3333
$pArr[62,0] = 5000
$pArr[63,0] = 4000
$pArr[64,0] = 4667
$pArr[65,0] = 4000
$pArr[66,0] = 5000
$pArr[67,0] = 4000
$pArr[68,0] = 4333
$pArr[69,0] = 3833
$pArr[70,0] = 3000
$pArr[71,0] = 5000
$pArr[72,0] = 4333
$pArr[73,0] = 4000
$pArr[74,0] = 3333
$pArr[75,0] = 3833
$pArr[76,0] = 3000
$pArr[77,0] = 4333
$pArr[78,0] = 3833
$pArr[79,0] = 3167
$pArr[80,0] = 5000
$pArr[81,0] = 4167
$pArr[82,0] = 3667
$pArr[83,0] = 5000
$pArr[84,0] = 4333
$pArr[85,0] = 5000
$pArr[86,0] = 4333
$pArr[87,0] = 3667
$pArr[88,0] = 4173
$pArr[89,0] = 3500
$pArr[90,0] = 2833
$pArr[91,0] = 3689
$pArr[92,0] = 2667
$pArr[93,0] = 4667
$pArr[94,0] = 4167
$pArr[95,0] = 3507
$pArr[96,0] = 4333
$pArr[97,0] = 3833
$pArr[98,0] = 3167
$pArr[99,0] = 4529
$pArr[100,0] = 4022
$pArr[101,0] = 3362
$pArr[102,0] = 5333
$pArr[103,0] = 3500
$pArr[104,0] = 4667
$pArr[105,0] = 3833
$pArr[106,0] = 3000
$pArr[107,0] = 3778
$pArr[108,0] = 2667
$pArr[109,0] = 4167
$pArr[110,0] = 3496
$pArr[111,0] = 2833
$pArr[112,0] = 4667
$pArr[113,0] = 3333
$pArr[114,0] = 4000
$pArr[115,0] = 3500
$pArr[116,0] = 3000
$pArr[117,0] = 4000
$pArr[118,0] = 3500
$pArr[119,0] = 2833
$ws.Range("P182:P301").Value = $pArr

# New rows 300-301 need the remaining (previously-constant) columns
# filled in, since they did not exist before.
$ws.Range("A300").Value = 9
$ws.Range("B300").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C300").Value = "Metropolitana"
$ws.Range("E300").Value = 13
$ws.Range("F300").Value = 100112009
$ws.Range("G300").Value = "Acelga"
$ws.Range("H300").Value = "Sin especificar"
$ws.Range("N300").Value = "`$/docena de atados"
$ws.Range("O300").Value = "Región Metropolitana"
$ws.Range("Q300").Value = 3
$ws.Range("R300").Value = "Hortaliza"
$ws.Range("D300").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A301").Value = 9
$ws.Range("B301").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C301").Value = "Metropolitana"
$ws.Range("E301").Value = 13
$ws.Range("F301").Value = 100112009
$ws.Range("G301").Value = "Acelga"
$ws.Range("H301").Value = "Sin especificar"
$ws.Range("N301").Value = "`$/docena de atados"
$ws.Range("O301").Value = "Región Metropolitana"
$ws.Range("Q301").Value = 3
$ws.Range("R301").Value = "Hortaliza"
$ws.Range("D301").NumberFormat = "YYYY-MM-DD HH:MM:SS"
